$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal-text value into a cell while preserving the
# cells original (default/no explicit format) style - prevents Excel
# from auto-coercing numeric-looking strings (e.g. "1.00", "607.10")
# into actual numbers, and avoids introducing a new style index.
function Set-TextValue {
    param($Range, [string]$Text)
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.Style = "Normal"
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "63.661.36"
Set-TextValue $ws.Range("E2") "  -4.06%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "3.094.80"
Set-TextValue $ws.Range("E3") "  -5.28%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.00"
Set-TextValue $ws.Range("E4") "  +0.06%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "607.10"
Set-TextValue $ws.Range("E5") "  -1.36%  "

# Row 6 - Solana
Set-TextValue $ws.Range("D6") "144.41"
Set-TextValue $ws.Range("E6") "  -8.45%  "

# Row 7 - USDC
Set-TextValue $ws.Range("E7") "  -0.02%  "

# Row 8 - LidoStakedEther
Set-TextValue $ws.Range("D8") "3.094.74"
Set-TextValue $ws.Range("E8") "  -5.23%  "

# Row 9 - XRP
Set-TextValue $ws.Range("E9") "  -4.90%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("E10") "  -7.55%  "

# Row 11 - Toncoin
Set-TextValue $ws.Range("D11") "5.22"
Set-TextValue $ws.Range("E11") "  -9.87%  "

# Row 12 - Cardano
Set-TextValue $ws.Range("D12") "0.467"
Set-TextValue $ws.Range("E12") "  -5.63%  "

# Row 13 - ShibaInu
Set-TextValue $ws.Range("E13") "  -8.62%  "

# Row 14 - Avalanche
Set-TextValue $ws.Range("D14") "35.08"
Set-TextValue $ws.Range("E14") "  -10.14%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue $ws.Range("D15") "3.598.44"
Set-TextValue $ws.Range("E15") "  -5.47%  "

# Row 16 - TRON
Set-TextValue $ws.Range("D16") "0.115"
Set-TextValue $ws.Range("E16") "  +0.83%  "

# Row 17 - WrappedBTC
Set-TextValue $ws.Range("D17") "63.674.19"
Set-TextValue $ws.Range("E17") "  -4.19%  "

# Row 18 - WrappedEther
Set-TextValue $ws.Range("D18") "3.092.90"
Set-TextValue $ws.Range("E18") "  -5.40%  "

# Row 19 - Polkadot
Set-TextValue $ws.Range("D19") "6.80"
Set-TextValue $ws.Range("E19") "  -8.53%  "

# Row 20 - BitcoinCash
Set-TextValue $ws.Range("D20") "474.06"
Set-TextValue $ws.Range("E20") "  -6.16%  "

# Row 21 - Chainlink
Set-TextValue $ws.Range("D21") "14.58"
Set-TextValue $ws.Range("E21") "  -6.08%  "

# Row 22 - Polygon
Set-TextValue $ws.Range("E22") "  -7.74%  "

# Row 23 - Uniswap
Set-TextValue $ws.Range("D23") "7.66"
Set-TextValue $ws.Range("E23") "  -5.82%  "

# Row 24 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D24") "13.50"
Set-TextValue $ws.Range("E24") "  -7.89%  "

# Row 25 - Litecoin
Set-TextValue $ws.Range("D25") "83.15"
Set-TextValue $ws.Range("E25") "  -4.46%  "

# Row 26 - Dai
Set-TextValue $ws.Range("D26") "0.999"
Set-TextValue $ws.Range("E26") "  -0.08%  "

# Row 27 - PancakeSwap
Set-TextValue $ws.Range("E27") "  -9.81%  "

# Row 28 - RenderToken
Set-TextValue $ws.Range("D28") "8.35"
Set-TextValue $ws.Range("E28") "  -9.52%  "

# Row 29 - ImmutableX
Set-TextValue $ws.Range("D29") "2.13"
Set-TextValue $ws.Range("E29") "  -11.30%  "

# Row 30 - NEARProtocol
Set-TextValue $ws.Range("E30") "  -5.33%  "

# Row 31 - Hedera
Set-TextValue $ws.Range("E31") "  -12.60%  "

# Row 32 - FirstDigitalUSD
Set-TextValue $ws.Range("E32") "  +0.02%  "

# Row 33 - Stacks
Set-TextValue $ws.Range("E33") "  -6.30%  "

# Row 34 - EthereumClassic
Set-TextValue $ws.Range("D34") "26.11"
Set-TextValue $ws.Range("E34") "  -6.83%  "

# Row 35 - Mantle
Set-TextValue $ws.Range("D35") "1.11"
Set-TextValue $ws.Range("E35") "  -3.95%  "

# Row 36 - Filecoin
Set-TextValue $ws.Range("D36") "5.91"
Set-TextValue $ws.Range("E36") "  -8.71%  "

# Row 37 - OKB
Set-TextValue $ws.Range("D37") "52.58"
Set-TextValue $ws.Range("E37") "  -5.48%  "

# Row 38 - PEPE
Set-TextValue $ws.Range("D38") "0.0₃0730"
Set-TextValue $ws.Range("E38") "  -7.80%  "

# Row 39 - Bittensor
Set-TextValue $ws.Range("D39") "457.96"
Set-TextValue $ws.Range("E39") "  -7.42%  "

# Row 40 - dogwifhat
Set-TextValue $ws.Range("D40") "2.92"
Set-TextValue $ws.Range("E40") "  -14.70%  "

# Row 41 - VeChain
Set-TextValue $ws.Range("E41") "  -8.27%  "

# Row 42 - Kaspa
Set-TextValue $ws.Range("E42") "  -8.11%  "

# Row 43 - Cosmos
Set-TextValue $ws.Range("D43") "8.30"
Set-TextValue $ws.Range("E43") "  -5.93%  "

# Row 44 - Maker
Set-TextValue $ws.Range("D44") "2.824.54"
Set-TextValue $ws.Range("E44") "  -6.29%  "

# Row 45 - TheGraph
Set-TextValue $ws.Range("E45") "  -9.53%  "

# Row 46 - Fetch.AI
Set-TextValue $ws.Range("D46") "2.24"
Set-TextValue $ws.Range("E46") "  -11.56%  "

# Row 47 - ThetaToken
Set-TextValue $ws.Range("D47") "2.40"
Set-TextValue $ws.Range("E47") "  -4.00%  "

# Row 49 - InjectiveProtocol
Set-TextValue $ws.Range("D49") "26.05"
Set-TextValue $ws.Range("E49") "  -10.09%  "

# Row 50 - Stellar
Set-TextValue $ws.Range("E50") "  -5.63%  "

# Row 51 - Monero
Set-TextValue $ws.Range("D51") "118.05"
Set-TextValue $ws.Range("E51") "  -2.53%  "
